# Add a "project name" column (aws-maven) in front of the existing
# CLASS / METHOD / COVERAGE columns of the cobertura coverage sheet.
#
# Before: A=CLASS        B=METHOD       C=COVERAGE
# After:  A=aws-maven(*) B=CLASS        C=METHOD        D=COVERAGE
#   (*) header row (row 1) keeps column A empty - only the data rows
#       (2..97) get the literal project name "aws-maven".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column A. This shifts the existing
# CLASS/METHOD/COVERAGE columns from A/B/C to B/C/D (data + column
# widths move together), exactly like the XML diff shows.
$ws.Columns("A").Insert()

# Populate the new column A with the project name for every data row
# (row 1 is the header row and is left blank in column A).
$ws.Range("A2:A97").Value2 = "aws-maven"

# Give the new column its own width (narrower than the others), matching
# the width added for column 1 in the diff's <cols> list. The engine
# quantizes COM ColumnWidth to the nearest 1/7 character, so this is the
# closest representable value to the target raw width (8.589887640449438).
$ws.Columns("A").ColumnWidth = 7.857142857142857
